$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 4 with new entry data
$ws.Range("A4").Value = "Tuesday 10.4.18"
$ws.Range("B4").Value = "2200-2300"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Completing Point and Polygon"

# Update selection to B24 as in the diff
$ws.Range("B24").Select()
